# Update countries & provincias Spain
# Refresh the COVID-19 country data table on sheet "Pais":
#  - bump the "Datos actualizados" timestamp
#  - a handful of rows swap which country label they show (alphabetical
#    re-sort of near-tied entries in the source feed) so both the label
#    (column A) and the numbers (columns B:H) for those rows change
#  - the remaining touched rows just get refreshed totals

function Set-RowValues {
    param($ws, $row, $startCol, $vals)
    $col = $startCol
    foreach ($v in $vals) {
        $ws.Cells.Item($row, $col).Value = $v
        $col = $col + 1
    }
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp banner (A1) -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 20 de Junio de 2020 a las 00:43"

# --- Plain data refreshes (country label unchanged) ------------------------
# Row 4: Estados Unidos
Set-RowValues $ws 4 2 @(2294246, 30595, 954624, 1218243, 0, 691, 121379)

# Row 53: Nigeria
Set-RowValues $ws 53 2 @(19147, 667, 6581, 12079, 0, 12, 487)

# Row 67: Chequia
Set-RowValues $ws 67 2 @(10406, 126, 7473, 2598)

# Row 69: Noruega
$ws.Cells.Item(69, 2).Value = 8726
$ws.Cells.Item(69, 3).Value = 18
$ws.Cells.Item(69, 5).Value = 344

# Row 76: Uzbekistan
$ws.Cells.Item(76, 2).Value = 5946
$ws.Cells.Item(76, 3).Value = 179
$ws.Cells.Item(76, 5).Value = 1654

# Row 89: Etiopia
$ws.Cells.Item(89, 4).Value = 1029
$ws.Cells.Item(89, 5).Value = 2969

# Row 90: Bulgaria
Set-RowValues $ws 90 2 @(3755, 81, 2008, 1554, 0, 3, 193)

# Row 135: Uruguay
$ws.Cells.Item(135, 2).Value = 853
$ws.Cells.Item(135, 3).Value = 3
$ws.Cells.Item(135, 5).Value = 15

# Row 165: Islas Caimanes
$ws.Cells.Item(165, 2).Value = 195
$ws.Cells.Item(165, 3).Value = 2
$ws.Cells.Item(165, 4).Value = 143

# --- Rows whose country label re-sorts with its neighbour -------------------
# Rows 160/161: Birmania <-> Surinam swap places, each keeps/gets its own
# refreshed totals.
$ws.Cells.Item(160, 1).Value = "Surinam"
Set-RowValues $ws 160 2 @(293, 16, 74, 211, 0, 1, 8)
$ws.Cells.Item(161, 1).Value = "Birmania"
Set-RowValues $ws 161 2 @(286, 23, 192, 88, 0, 0, 6)

# Rows 170/171: Angola <-> Guadalupe swap places.
$ws.Cells.Item(170, 1).Value = "Guadalupe"
Set-RowValues $ws 170 2 @(174, 3, 157, 3, 0, 0, 14)
$ws.Cells.Item(171, 1).Value = "Angola"
Set-RowValues $ws 171 2 @(172, 6, 66, 98, 0, 0, 8)

# Rows 202/203: Dominica <-> Fiyi swap places (identical totals, so only the
# label moves).
$ws.Cells.Item(202, 1).Value = "Fiyi"
$ws.Cells.Item(203, 1).Value = "Dominica"

# Rows 208/209: Islas Turcas y Caicos <-> Santa Sede swap places.
$ws.Cells.Item(208, 1).Value = "Santa Sede"
Set-RowValues $ws 208 2 @(12, 0, 12, 0, 0, 0, 0)
$ws.Cells.Item(209, 1).Value = "Islas Turcas y Caicos"
Set-RowValues $ws 209 2 @(12, 0, 11, 0, 0, 0, 1)
